$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force text format on D:E so numeric-looking strings
# (e.g. "1.000", "0.9995") keep their exact text representation
# instead of being auto-converted to numbers by Excel input parsing.
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '31.039.22'
$ws.Range('E2').Value = '  +1.30%  '
$ws.Range('D3').Value = '1.959.59'
$ws.Range('E3').Value = '  -0.16%  '
$ws.Range('D4').Value = '0.9995'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '245.82'
$ws.Range('E5').Value = '  -1.35%  '
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('D7').Value = '0.4884'
$ws.Range('E7').Value = '  +0.52%  '
$ws.Range('D8').Value = '0.2967'
$ws.Range('E8').Value = '  +0.69%  '
$ws.Range('D9').Value = '0.06835'
$ws.Range('E9').Value = '  +0.80%  '
$ws.Range('D10').Value = '19.23'
$ws.Range('E10').Value = '  -1.35%  '
$ws.Range('D11').Value = '107.01'
$ws.Range('E11').Value = '  -2.19%  '
$ws.Range('D12').Value = '1.960.36'
$ws.Range('E12').Value = '  -0.10%  '
$ws.Range('D13').Value = '0.07837'
$ws.Range('E13').Value = '  +0.87%  '
$ws.Range('E14').Value = '  +0.94%  '
$ws.Range('D15').Value = '0.7055'
$ws.Range('E15').Value = '  +2.67%  '
$ws.Range('D16').Value = '284.52'
$ws.Range('E16').Value = '  -3.24%  '
$ws.Range('D17').Value = '31.063.56'
$ws.Range('E17').Value = '  +1.35%  '
$ws.Range('E18').Value = '  +0.15%  '
$ws.Range('D19').Value = '0.000007716'
$ws.Range('E19').Value = '  +0.23%  '
$ws.Range('D20').Value = '2.224.26'
$ws.Range('E20').Value = '  +0.06%  '
$ws.Range('D21').Value = '1.000'
$ws.Range('E21').Value = '  +0.10%  '
$ws.Range('E22').Value = '  -1.63%  '
$ws.Range('D23').Value = '0.9998'
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').Value = '6.523'
$ws.Range('E24').Value = '  -1.55%  '
$ws.Range('D25').Value = '9.834'
$ws.Range('E25').Value = '  -0.43%  '
$ws.Range('D26').Value = '169.12'
$ws.Range('E26').Value = '  -0.70%  '
$ws.Range('D27').Value = '20.05'
$ws.Range('E27').Value = '  -0.37%  '
$ws.Range('D28').Value = '2.210'
$ws.Range('E28').Value = '  +1.61%  '
$ws.Range('D29').Value = '0.1058'
$ws.Range('E29').Value = '  -0.72%  '
$ws.Range('D30').Value = '1.396'
$ws.Range('E30').Value = '  -2.88%  '
$ws.Range('D31').Value = '1.586'
$ws.Range('E31').Value = '  -1.38%  '
$ws.Range('D32').Value = '4.611'
$ws.Range('E32').Value = '  -1.74%  '
$ws.Range('D33').Value = '4.453'
$ws.Range('E33').Value = '  +0.17%  '
$ws.Range('D34').Value = '0.04952'
$ws.Range('D35').Value = '0.7643'
$ws.Range('E35').Value = '  -0.86%  '
$ws.Range('D36').Value = '1.176'
$ws.Range('E36').Value = '  -0.49%  '
$ws.Range('D37').Value = '2.733'
$ws.Range('E37').Value = '  -0.04%  '
$ws.Range('E38').Value = '  -0.95%  '
$ws.Range('D39').Value = '2.702'
$ws.Range('E39').Value = '  -0.69%  '
$ws.Range('D40').Value = '6.567'
$ws.Range('E40').Value = '  +1.40%  '
$ws.Range('D41').Value = '78.39'
$ws.Range('E41').Value = '  +11.45%  '
$ws.Range('D42').Value = '2.120'
$ws.Range('E42').Value = '  -0.36%  '
$ws.Range('D43').Value = '0.9093'
$ws.Range('E43').Value = '  +3.85%  '
$ws.Range('D45').Value = '109.19'
$ws.Range('E45').Value = '  -0.06%  '
$ws.Range('D46').Value = '8.184'
$ws.Range('E46').Value = '  +8.98%  '
$ws.Range('D47').Value = '1.001'
$ws.Range('E47').Value = '  +0.01%  '
$ws.Range('D48').Value = '1.022.89'
$ws.Range('E48').Value = '  +11.68%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '9.414'
$ws.Range('E49').Value = '  +0.66%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').Value = '0.1264'
$ws.Range('E50').Value = '  -1.59%  '
$ws.Range('D51').Value = '36.01'
$ws.Range('E51').Value = '  -0.13%  '

# Restore the original (default) cell style now that text is committed,
# so no stray number-format style lingers on these cells.
$priceRange.Style = "Normal"
